# test precision for EMU->pt->EMU round trip through COM Height
$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$s = $m.Shapes.Item(2)  # sldNum placeholder
$targetEmu = 270360
$pt = $targetEmu / 12700.0
Write-Host ("pt=" + $pt)
$s.Height = $pt
Write-Host ("Height after set=" + $s.Height)
